$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'290.53"
$ws.Range("E2").Value = "'-3.84%"
$ws.Range("G2").Value = "'7"
$ws.Range("D3").Value = "'30.73"
$ws.Range("E3").Value = "'-6.00%"
$ws.Range("G3").Value = "'7"
$ws.Range("E4").Value = "'-0.23%"
$ws.Range("G4").Value = "'7"
$ws.Range("D5").Value = "'0.07217"
$ws.Range("E5").Value = "'-7.46%"
$ws.Range("G5").Value = "'7"
$ws.Range("D6").Value = "'1.815"
$ws.Range("E6").Value = "'-8.51%"
$ws.Range("G6").Value = "'7"
$ws.Range("D7").Value = "'7.675"
$ws.Range("E7").Value = "'-2.35%"
$ws.Range("G7").Value = "'7"
$ws.Range("D8").Value = "'3.725"
$ws.Range("E8").Value = "'-1.95%"
$ws.Range("G8").Value = "'7"
$ws.Range("D9").Value = "'0.8964"
$ws.Range("E9").Value = "'-3.31%"
$ws.Range("G9").Value = "'7"
$ws.Range("D10").Value = "'0.1649"
$ws.Range("E10").Value = "'-6.36%"
$ws.Range("G10").Value = "'7"
$ws.Range("D11").Value = "'0.07649"
$ws.Range("E11").Value = "'-2.94%"
$ws.Range("G11").Value = "'7"
$ws.Range("D12").Value = "'0.07970"
$ws.Range("E12").Value = "'-7.61%"
$ws.Range("G12").Value = "'7"
$ws.Range("D13").Value = "'0.03038"
$ws.Range("E13").Value = "'-4.10%"
$ws.Range("G13").Value = "'7"
$ws.Range("E14").Value = "'0.06%"
$ws.Range("G14").Value = "'7"
$ws.Range("D15").Value = "'0.001508"
$ws.Range("E15").Value = "'-1.36%"
$ws.Range("G15").Value = "'7"
$ws.Range("D16").Value = "'0.005709"
$ws.Range("E16").Value = "'-0.51%"
$ws.Range("G16").Value = "'7"
$ws.Range("D17").Value = "'3.460"
$ws.Range("E17").Value = "'-0.09%"
$ws.Range("G17").Value = "'7"
$ws.Range("D18").Value = "'2.083"
$ws.Range("E18").Value = "'-3.30%"
$ws.Range("G18").Value = "'7"
$ws.Range("E19").Value = "'-0.42%"
$ws.Range("G19").Value = "'7"
$ws.Range("D20").Value = "'0.1331"
$ws.Range("E20").Value = "'0.90%"
$ws.Range("G20").Value = "'7"
$ws.Range("D21").Value = "'4.035"
$ws.Range("E21").Value = "'-6.92%"
$ws.Range("G21").Value = "'7"
$ws.Range("E22").Value = "'5.36%"
$ws.Range("G22").Value = "'7"
$ws.Range("D23").Value = "'0.04519"
$ws.Range("E23").Value = "'-0.99%"
$ws.Range("G23").Value = "'7"
$ws.Range("E24").Value = "'-1.02%"
$ws.Range("G24").Value = "'7"
$ws.Range("D25").Value = "'0.004016"
$ws.Range("E25").Value = "'-9.52%"
$ws.Range("G25").Value = "'7"
$ws.Range("D26").Value = "'0.0001251"
$ws.Range("E26").Value = "'-0.09%"
$ws.Range("G26").Value = "'7"
$ws.Range("G27").Value = "'7"
$ws.Range("G28").Value = "'7"
$ws.Range("G29").Value = "'7"
$ws.Range("G30").Value = "'7"
$ws.Range("G31").Value = "'7"
$ws.Range("G32").Value = "'7"
$ws.Range("G33").Value = "'7"
$ws.Range("G34").Value = "'7"
$ws.Range("G35").Value = "'7"
$ws.Range("G36").Value = "'7"
$ws.Range("G37").Value = "'7"
$ws.Range("G38").Value = "'7"
$ws.Range("D39").Value = "'0.01604"
$ws.Range("E39").Value = "'-6.01%"
$ws.Range("G39").Value = "'7"
$ws.Range("D40").Value = "'0.04395"
$ws.Range("E40").Value = "'-7.12%"
$ws.Range("G40").Value = "'7"
$ws.Range("D41").Value = "'0.007377"
$ws.Range("E41").Value = "'-6.01%"
$ws.Range("G41").Value = "'7"
$ws.Range("D42").Value = "'0.1310"
$ws.Range("E42").Value = "'-3.33%"
$ws.Range("G42").Value = "'7"
$ws.Range("D43").Value = "'0.007701"
$ws.Range("G43").Value = "'7"
$ws.Range("D44").Value = "'0.002051"
$ws.Range("E44").Value = "'-12.47%"
$ws.Range("G44").Value = "'7"
$ws.Range("D45").Value = "'0.009214"
$ws.Range("E45").Value = "'-12.35%"
$ws.Range("G45").Value = "'7"
$ws.Range("D46").Value = "'0.00005882"
$ws.Range("E46").Value = "'-6.19%"
$ws.Range("G46").Value = "'7"
$ws.Range("E47").Value = "'-0.09%"
$ws.Range("G47").Value = "'7"
$ws.Range("D48").Value = "'2.247"
$ws.Range("E48").Value = "'173.92%"
$ws.Range("G48").Value = "'7"
$ws.Range("E49").Value = "'-3.37%"
$ws.Range("G49").Value = "'7"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.09%"
$ws.Range("G50").Value = "'7"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.09%"
$ws.Range("G51").Value = "'7"
